# "Generate Report for Handoff"
#
# A new handoff xliff-generation pass ran for the four files that were
# still sitting at "Ready for handoff" (3cc9f465, a6145d5c, c5fcb705,
# fd546b82). For those rows, on both the zh-cn and de-de localization
# sheets:
#   - Priority flips from "low" to "ht" (high throughput / handoff-triggered)
#   - The Latest Handoff Datetime is refreshed to the moment this handoff
#     report was generated (zh-cn: 06:36:16 -> 06:36:32,
#     de-de: 06:36:22 -> 06:36:36)
#
# The de-de handoff timestamp is the same value shown on the Overview
# sheet's "Latest HO Xliff Generate Date" column for those rows, so that
# gets refreshed too.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$rows = 4, 5, 6, 7

foreach ($r in $rows) {
    # Overview: Latest HO Xliff Generate Date (column G)
    $overview.Range("G$r").Value = "2016-09-06 06:36:36"

    # zh-cn: Priority (E) + Latest Handoff Datetime (H)
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-09-06 06:36:32"

    # de-de: Priority (E) + Latest Handoff Datetime (H)
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-09-06 06:36:36"
}

$wb.Save()
